$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "34.068.82"
$ws.Range("E2").Value = "  +9.98%  "
$ws.Range("D3").Value = "1.788.22"
$ws.Range("E3").Value = "  +6.52%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.997"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.16%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "225.17"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.41%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.558"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +4.67%  "
$ws.Range("E7").Value = "  -0.16%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "30.93"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +5.97%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "46.31"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +3.34%  "
$ws.Range("E10").Value = "  +5.25%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0661"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +3.43%  "
$ws.Range("E12").Value = "  +1.68%  "
$ws.Range("E13").Value = "  +6.60%  "
$ws.Range("D14").Value = "1.787.02"
$ws.Range("E14").Value = "  +6.69%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.628"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +3.82%  "
$ws.Range("D16").Value = "33.997.09"
$ws.Range("E16").Value = "  +9.81%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "9.99"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.31%  "
$ws.Range("E18").Value = "  +2.38%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "68.55"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +3.51%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "251.99"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.89%  "
$ws.Range("D21").Value = "0.0₃0740"
$ws.Range("E21").Value = "  +2.87%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.995"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.35%  "
$ws.Range("E23").Value = "  +2.98%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "4.22"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.87%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.16"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.10%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "157.53"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.83%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "16.47"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +3.78%  "
$ws.Range("E28").Value = "  +2.09%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.94"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +3.73%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.998"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.14%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.80"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +8.55%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0511"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +3.29%  "
$ws.Range("E33").Value = "  +4.39%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.53"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +6.00%  "
$ws.Range("D35").Value = "1.495.09"
$ws.Range("E35").Value = "  -1.45%  "
$ws.Range("E36").Value = "  +1.49%  "
$ws.Range("E37").Value = "  +3.73%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.629"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +3.02%  "
$ws.Range("E39").Value = "  +3.23%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "83.09"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.70%  "
$ws.Range("E41").Value = "  +2.76%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.71"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.31%  "
$ws.Range("E43").Value = "  +6.07%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.08"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.18%  "
$ws.Range("E45").Value = "  +1.24%  "
$ws.Range("E46").Value = "  +3.68%  "
$ws.Range("E47").Value = "  +7.15%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "5.75"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.99%  "
$ws.Range("E49").Value = "  +0.00%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "11.96"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +12.79%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "50.94"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.78%  "
